$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F column (time_taken) timestamps on the "data" sheet (rows 2-56) ---
$timestamps = @(
    "2021-10-05 14:22:05.997916",
    "2021-10-05 14:22:05.997924",
    "2021-10-05 14:22:05.997928",
    "2021-10-05 14:22:05.997930",
    "2021-10-05 14:22:05.997933",
    "2021-10-05 14:22:05.997936",
    "2021-10-05 14:22:05.997939",
    "2021-10-05 14:22:05.997946",
    "2021-10-05 14:22:05.997949",
    "2021-10-05 14:22:05.997952",
    "2021-10-05 14:22:05.997955",
    "2021-10-05 14:22:05.997958",
    "2021-10-05 14:22:05.997960",
    "2021-10-05 14:22:05.997963",
    "2021-10-05 14:22:05.997966",
    "2021-10-05 14:22:05.997969",
    "2021-10-05 14:22:05.997972",
    "2021-10-05 14:22:05.997974",
    "2021-10-05 14:22:05.997977",
    "2021-10-05 14:22:05.997980",
    "2021-10-05 14:22:05.997983",
    "2021-10-05 14:22:05.997986",
    "2021-10-05 14:22:05.997988",
    "2021-10-05 14:22:05.997991",
    "2021-10-05 14:22:05.997994",
    "2021-10-05 14:22:05.997997",
    "2021-10-05 14:22:05.998000",
    "2021-10-05 14:22:05.998002",
    "2021-10-05 14:22:05.998005",
    "2021-10-05 14:22:05.998008",
    "2021-10-05 14:22:05.998010",
    "2021-10-05 14:22:05.998013",
    "2021-10-05 14:22:05.998016",
    "2021-10-05 14:22:05.998019",
    "2021-10-05 14:22:05.998022",
    "2021-10-05 14:22:05.998025",
    "2021-10-05 14:22:05.998027",
    "2021-10-05 14:22:05.998030",
    "2021-10-05 14:22:05.998033",
    "2021-10-05 14:22:05.998035",
    "2021-10-05 14:22:05.998038",
    "2021-10-05 14:22:05.998041",
    "2021-10-05 14:22:05.998044",
    "2021-10-05 14:22:05.998047",
    "2021-10-05 14:22:05.998050",
    "2021-10-05 14:22:05.998052",
    "2021-10-05 14:22:05.998055",
    "2021-10-05 14:22:05.998058",
    "2021-10-05 14:22:05.998060",
    "2021-10-05 14:22:05.998063",
    "2021-10-05 14:22:05.998066",
    "2021-10-05 14:22:05.998068",
    "2021-10-05 14:22:05.998071",
    "2021-10-05 14:22:05.998074",
    "2021-10-05 14:22:05.998077"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $dataSheet.Cells.Item($i + 2, 6).Value = $timestamps[$i]
}

# --- Add a new "metadata" worksheet positioned right after the "data" sheet ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1)
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row (A2:G2)
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Paediatric disorders - additional genes"
$metaSheet.Cells.Item(2, 3).Value = 479
# data_version ("1.94") must stay a text string, not be coerced to a number
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.94"
$metaSheet.Cells.Item(2, 5).Value = "2021-07-08T10:47:41.226728Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:22:05.994648"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/479/?format=json"

# Match the bold/centered/bordered header style used on the "data" sheet (style index 1),
# and the same style used for the "data" sheet's A-column index cells (A2:A56).
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$dataSheet.Select()
